# Generate Report for Handback
# Marks the 0d1ca4e8-... file as handed back (in sync with en-US) for both
# the zh-cn and de-de locales, recording the target file / handback file /
# handback datetime for each, and reflects the new status on the Overview
# sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: row 2 is the 0d1ca4e8-... file; its zh-cn / de-de
# status columns now report the handback.
# ---------------------------------------------------------------------
$overview.Range("B2").Value = $statusHandedBack
$overview.Range("C2").Value = $statusHandedBack

# ---------------------------------------------------------------------
# zh-cn sheet: row 2 is the 0d1ca4e8-... file.
# ---------------------------------------------------------------------
$zhcn.Range("C2").Value = $statusHandedBack
$zhcn.Range("F2").Value = "0d1ca4e8-563f-4906-8bcc-a3977a07398f.md"
$zhcn.Range("F2").Style = "HyperLink"
$zhcn.Range("G2").Value = "0d1ca4e8-563f-4906-8bcc-a3977a07398f.ccb5949cadfba9ae28124f850d36e8217cf49b07.zh-cn.xlf"
$zhcn.Range("G2").Style = "HyperLink"
$zhcn.Range("H2").Value = "2016-03-19 06:29:27"

$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/129206683dae6d6731884191686087d63b4502f2/e2e/0d1ca4e8-563f-4906-8bcc-a3977a07398f.md", "", "", "0d1ca4e8-563f-4906-8bcc-a3977a07398f.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e56cbefbddcddda91ef90bb908aef1f844b14188/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0d1ca4e8-563f-4906-8bcc-a3977a07398f.ccb5949cadfba9ae28124f850d36e8217cf49b07.zh-cn.xlf", "", "", "0d1ca4e8-563f-4906-8bcc-a3977a07398f.ccb5949cadfba9ae28124f850d36e8217cf49b07.zh-cn.xlf")

# ---------------------------------------------------------------------
# de-de sheet: row 2 is the 0d1ca4e8-... file.
# ---------------------------------------------------------------------
$dede.Range("C2").Value = $statusHandedBack
$dede.Range("F2").Value = "0d1ca4e8-563f-4906-8bcc-a3977a07398f.md"
$dede.Range("F2").Style = "HyperLink"
$dede.Range("G2").Value = "0d1ca4e8-563f-4906-8bcc-a3977a07398f.ccb5949cadfba9ae28124f850d36e8217cf49b07.de-de.xlf"
$dede.Range("G2").Style = "HyperLink"
$dede.Range("H2").Value = "2016-03-19 06:29:32"

$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/129206683dae6d6731884191686087d63b4502f2/e2e/0d1ca4e8-563f-4906-8bcc-a3977a07398f.md", "", "", "0d1ca4e8-563f-4906-8bcc-a3977a07398f.md")
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dc56cd0ebbfc96641b3435d82f6bfb21aabe30b4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0d1ca4e8-563f-4906-8bcc-a3977a07398f.ccb5949cadfba9ae28124f850d36e8217cf49b07.de-de.xlf", "", "", "0d1ca4e8-563f-4906-8bcc-a3977a07398f.ccb5949cadfba9ae28124f850d36e8217cf49b07.de-de.xlf")
